$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:A7").Value = 900

$ws.Range("A8").Select()
